$d = $word.ActiveDocument

# Paragraph 1: "Hooverville" heading -> change style to Title, remove bookmark
$p1 = $d.Paragraphs.Item(1)
$p1.Range.Style = "Title"

# Paragraph 2: "By Dorothy Day" (bold) -> replace with Authors-styled paragraph
# containing three runs: "Dorothy", " ", "Day" (no bold)
$p2 = $d.Paragraphs.Item(2)
$p2.Range.Text = "Dorothy Day"
$p2.Range.Style = "Authors"
$p2.Range.Font.Bold = 0

$d.Bookmarks.Item("hooverville").Delete()
